$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Find.Execute("2023-06-06 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-07 Wednesday", 2) | Out-Null

# Update each of the 100 table cells (20 rows x 5 columns), in row-major order
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "75-39="
$t.Cell(1,2).Range.Text = "68-5="
$t.Cell(1,3).Range.Text = "63-18="
$t.Cell(1,4).Range.Text = "76-5="
$t.Cell(1,5).Range.Text = "40+11="
$t.Cell(2,1).Range.Text = "78-9="
$t.Cell(2,2).Range.Text = "31+59="
$t.Cell(2,3).Range.Text = "35-24="
$t.Cell(2,4).Range.Text = "94-67="
$t.Cell(2,5).Range.Text = "64-23="
$t.Cell(3,1).Range.Text = "36+1="
$t.Cell(3,2).Range.Text = "51-5="
$t.Cell(3,3).Range.Text = "1+6="
$t.Cell(3,4).Range.Text = "8+70="
$t.Cell(3,5).Range.Text = "23+23="
$t.Cell(4,1).Range.Text = "41+22="
$t.Cell(4,2).Range.Text = "17-6="
$t.Cell(4,3).Range.Text = "46-38="
$t.Cell(4,4).Range.Text = "3+49="
$t.Cell(4,5).Range.Text = "65-61="
$t.Cell(5,1).Range.Text = "97-13="
$t.Cell(5,2).Range.Text = "33+5="
$t.Cell(5,3).Range.Text = "88-46="
$t.Cell(5,4).Range.Text = "28+12="
$t.Cell(5,5).Range.Text = "47-47="
$t.Cell(6,1).Range.Text = "14+57="
$t.Cell(6,2).Range.Text = "19+41="
$t.Cell(6,3).Range.Text = "61-10="
$t.Cell(6,4).Range.Text = "61-1="
$t.Cell(6,5).Range.Text = "87-8="
$t.Cell(7,1).Range.Text = "11+81="
$t.Cell(7,2).Range.Text = "15+64="
$t.Cell(7,3).Range.Text = "20+31="
$t.Cell(7,4).Range.Text = "93-84="
$t.Cell(7,5).Range.Text = "23-0="
$t.Cell(8,1).Range.Text = "57-56="
$t.Cell(8,2).Range.Text = "2+4="
$t.Cell(8,3).Range.Text = "15+51="
$t.Cell(8,4).Range.Text = "68-1="
$t.Cell(8,5).Range.Text = "99-98="
$t.Cell(9,1).Range.Text = "43-28="
$t.Cell(9,2).Range.Text = "0+20="
$t.Cell(9,3).Range.Text = "65-20="
$t.Cell(9,4).Range.Text = "26+57="
$t.Cell(9,5).Range.Text = "38+32="
$t.Cell(10,1).Range.Text = "38-5="
$t.Cell(10,2).Range.Text = "60-11="
$t.Cell(10,3).Range.Text = "6+17="
$t.Cell(10,4).Range.Text = "4+62="
$t.Cell(10,5).Range.Text = "5+85="
$t.Cell(11,1).Range.Text = "22+6="
$t.Cell(11,2).Range.Text = "23+0="
$t.Cell(11,3).Range.Text = "68+21="
$t.Cell(11,4).Range.Text = "49-18="
$t.Cell(11,5).Range.Text = "65+8="
$t.Cell(12,1).Range.Text = "39-4="
$t.Cell(12,2).Range.Text = "67-53="
$t.Cell(12,3).Range.Text = "84-82="
$t.Cell(12,4).Range.Text = "65-25="
$t.Cell(12,5).Range.Text = "42+50="
$t.Cell(13,1).Range.Text = "29-17="
$t.Cell(13,2).Range.Text = "29+5="
$t.Cell(13,3).Range.Text = "98-77="
$t.Cell(13,4).Range.Text = "29+46="
$t.Cell(13,5).Range.Text = "69-45="
$t.Cell(14,1).Range.Text = "88-76="
$t.Cell(14,2).Range.Text = "28+34="
$t.Cell(14,3).Range.Text = "37+48="
$t.Cell(14,4).Range.Text = "32-25="
$t.Cell(14,5).Range.Text = "94-57="
$t.Cell(15,1).Range.Text = "22-15="
$t.Cell(15,2).Range.Text = "24+29="
$t.Cell(15,3).Range.Text = "64-37="
$t.Cell(15,4).Range.Text = "86-53="
$t.Cell(15,5).Range.Text = "69-49="
$t.Cell(16,1).Range.Text = "99-20="
$t.Cell(16,2).Range.Text = "8+2="
$t.Cell(16,3).Range.Text = "73-29="
$t.Cell(16,4).Range.Text = "98-88="
$t.Cell(16,5).Range.Text = "88-50="
$t.Cell(17,1).Range.Text = "3+78="
$t.Cell(17,2).Range.Text = "13+25="
$t.Cell(17,3).Range.Text = "30-5="
$t.Cell(17,4).Range.Text = "21+44="
$t.Cell(17,5).Range.Text = "80-16="
$t.Cell(18,1).Range.Text = "87-77="
$t.Cell(18,2).Range.Text = "6+65="
$t.Cell(18,3).Range.Text = "24+58="
$t.Cell(18,4).Range.Text = "18+37="
$t.Cell(18,5).Range.Text = "38+27="
$t.Cell(19,1).Range.Text = "52+24="
$t.Cell(19,2).Range.Text = "97-61="
$t.Cell(19,3).Range.Text = "56+8="
$t.Cell(19,4).Range.Text = "67-14="
$t.Cell(19,5).Range.Text = "35+7="
$t.Cell(20,1).Range.Text = "34+43="
$t.Cell(20,2).Range.Text = "86-0="
$t.Cell(20,3).Range.Text = "83-56="
$t.Cell(20,4).Range.Text = "55+3="
$t.Cell(20,5).Range.Text = "58+16="
